$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update quarterly financial figures (Income Statement, Balance Sheet, Cash Flow Statement)
# Row 8
$ws.Range("D8").Value = 4781800
$ws.Range("E8").Value = 4415100
$ws.Range("F8").Value = 3401900
$ws.Range("G8").Value = 3070900
$ws.Range("H8").Value = 1867100
$ws.Range("I8").Value = 2710000
$ws.Range("J8").Value = 2692800

# Row 9
$ws.Range("D9").Value = 2948700
$ws.Range("E9").Value = 2808700
$ws.Range("F9").Value = 2201100
$ws.Range("G9").Value = 2101400
$ws.Range("H9").Value = 1432000
$ws.Range("I9").Value = 2274900
$ws.Range("J9").Value = 2214500

# Row 10
$ws.Range("D10").Value = 1833100
$ws.Range("E10").Value = 1606400
$ws.Range("F10").Value = 1200800
$ws.Range("G10").Value = 969400
$ws.Range("H10").Value = 435100
$ws.Range("I10").Value = 435100
$ws.Range("J10").Value = 478200

# Row 15
$ws.Range("J15").Value = 1700

# Row 17
$ws.Range("D17").Value = 3734200
$ws.Range("E17").Value = 3703100
$ws.Range("F17").Value = 2680900
$ws.Range("G17").Value = 2799700
$ws.Range("H17").Value = 1889600
$ws.Range("I17").Value = 2698200
$ws.Range("J17").Value = 2636700

# Row 18
$ws.Range("D18").Value = 1047600
$ws.Range("E18").Value = 712000
$ws.Range("F18").Value = 721100
$ws.Range("G18").Value = 271200
$ws.Range("H18").Value = -22500
$ws.Range("I18").Value = 11800
$ws.Range("J18").Value = 56000

# Row 20
$ws.Range("D20").Value = 474700
$ws.Range("E20").Value = 499200
$ws.Range("F20").Value = 224700
$ws.Range("G20").Value = 268200
$ws.Range("H20").Value = 254400
$ws.Range("I20").Value = 190100
$ws.Range("J20").Value = 203200

# Row 21
$ws.Range("D21").Value = 1997200
$ws.Range("E21").Value = 1224500
$ws.Range("F21").Value = 1284300
$ws.Range("G21").Value = 561800
$ws.Range("H21").Value = 529800
$ws.Range("I21").Value = 232900
$ws.Range("J21").Value = 574800

# Row 22
$ws.Range("D22").Value = 269400
$ws.Range("E22").Value = 286800
$ws.Range("F22").Value = 196300
$ws.Range("G22").Value = 203100
$ws.Range("H22").Value = 168100
$ws.Range("I22").Value = 145700
$ws.Range("J22").Value = 223000

# Row 23
$ws.Range("D23").Value = 1252800
$ws.Range("E23").Value = 924400
$ws.Range("F23").Value = 749400
$ws.Range("G23").Value = 336300
$ws.Range("H23").Value = 63700
$ws.Range("I23").Value = 56100
$ws.Range("J23").Value = 36200

# Row 24
$ws.Range("D24").Value = 370100
$ws.Range("E24").Value = 235100
$ws.Range("F24").Value = 167800
$ws.Range("G24").Value = 110500
$ws.Range("H24").Value = 10800
$ws.Range("I24").Value = 36700
$ws.Range("J24").Value = 36000

# Row 26
$ws.Range("D26").Value = 882700
$ws.Range("E26").Value = 689300
$ws.Range("F26").Value = 581600
$ws.Range("G26").Value = 225800
$ws.Range("H26").Value = 52900
$ws.Range("I26").Value = 19500

# Row 27
$ws.Range("D27").Value = 686100
$ws.Range("E27").Value = 580500
$ws.Range("F27").Value = 512200
$ws.Range("G27").Value = 189100
$ws.Range("H27").Value = 55700
$ws.Range("I27").Value = 31900
$ws.Range("J27").Value = -7500

# Row 32
$ws.Range("D32").Value = -474700
$ws.Range("E32").Value = -499200
$ws.Range("F32").Value = -224700
$ws.Range("G32").Value = -268200
$ws.Range("H32").Value = -254400
$ws.Range("I32").Value = -190100
$ws.Range("J32").Value = -203200

# Row 33
$ws.Range("D33").Value = 686100
$ws.Range("E33").Value = 580500
$ws.Range("F33").Value = 512200
$ws.Range("G33").Value = 189100
$ws.Range("H33").Value = 55700
$ws.Range("I33").Value = 31900
$ws.Range("J33").Value = -7500

# Row 35
$ws.Range("D35").Value = 686100
$ws.Range("E35").Value = 580500
$ws.Range("F35").Value = 512200
$ws.Range("G35").Value = 189100
$ws.Range("H35").Value = 55700
$ws.Range("I35").Value = 31900
$ws.Range("J35").Value = -7500

# Row 41
$ws.Range("D41").Value = 3524300
$ws.Range("E41").Value = 3127500
$ws.Range("F41").Value = 2487300
$ws.Range("G41").Value = 2437300
$ws.Range("H41").Value = 2206300
$ws.Range("I41").Value = 2994200
$ws.Range("J41").Value = 2240500

# Row 42
$ws.Range("E42").Value = 10300
$ws.Range("F42").Value = 20200
$ws.Range("G42").Value = 19400

# Row 43
$ws.Range("D43").Value = 2972000
$ws.Range("E43").Value = 3122800
$ws.Range("F43").Value = 1925700
$ws.Range("G43").Value = 1761500
$ws.Range("H43").Value = 1353100
$ws.Range("I43").Value = 1179600
$ws.Range("J43").Value = 1245100

# Row 44
$ws.Range("D44").Value = 1283000
$ws.Range("E44").Value = 891600
$ws.Range("F44").Value = 896000
$ws.Range("G44").Value = 750800
$ws.Range("H44").Value = 919400
$ws.Range("I44").Value = 676700
$ws.Range("J44").Value = 737200

# Row 45
$ws.Range("D45").Value = 1655300
$ws.Range("E45").Value = 2531400
$ws.Range("F45").Value = 2115400
$ws.Range("G45").Value = 1191800
$ws.Range("H45").Value = 1011900
$ws.Range("I45").Value = 2393600
$ws.Range("J45").Value = 1756400

# Row 46
$ws.Range("D46").Value = 9434700
$ws.Range("E46").Value = 9683500
$ws.Range("F46").Value = 7444600
$ws.Range("G46").Value = 6160700
$ws.Range("H46").Value = 5490800
$ws.Range("I46").Value = 7244100
$ws.Range("J46").Value = 5979300

# Row 47
$ws.Range("D47").Value = 3191000
$ws.Range("E47").Value = 2826000
$ws.Range("F47").Value = 2111700
$ws.Range("G47").Value = 2006500
$ws.Range("H47").Value = 2264200
$ws.Range("I47").Value = 817300
$ws.Range("J47").Value = 918700

# Row 48
$ws.Range("D48").Value = 7747800
$ws.Range("E48").Value = 7671300
$ws.Range("F48").Value = 6422200
$ws.Range("G48").Value = 6226200
$ws.Range("H48").Value = 6986400
$ws.Range("I48").Value = 6769900
$ws.Range("J48").Value = 6541800

# Row 49
$ws.Range("D49").Value = 7515500
$ws.Range("E49").Value = 7700400
$ws.Range("F49").Value = 6311100
$ws.Range("G49").Value = 6270800
$ws.Range("H49").Value = 5195600
$ws.Range("I49").Value = 5126000
$ws.Range("J49").Value = 5637000

# Row 52
$ws.Range("D52").Value = 1241000
$ws.Range("E52").Value = 1401900
$ws.Range("F52").Value = 1427100
$ws.Range("G52").Value = 1219500
$ws.Range("H52").Value = 1196100
$ws.Range("I52").Value = 1187000
$ws.Range("J52").Value = 868600

# Row 54
$ws.Range("D54").Value = 29130000
$ws.Range("E54").Value = 29283200
$ws.Range("F54").Value = 23716700
$ws.Range("G54").Value = 21883900
$ws.Range("H54").Value = 21133000
$ws.Range("I54").Value = 21144300
$ws.Range("J54").Value = 19945400

# Row 57
$ws.Range("D57").Value = 4548100
$ws.Range("E57").Value = 4642700
$ws.Range("F57").Value = 2822300
$ws.Range("G57").Value = 2483300
$ws.Range("H57").Value = 1900800
$ws.Range("I57").Value = 1989700
$ws.Range("J57").Value = 1663700

# Row 58
$ws.Range("D58").Value = 3962500
$ws.Range("E58").Value = 4118600
$ws.Range("F58").Value = 4768000
$ws.Range("G58").Value = 4562400
$ws.Range("H58").Value = 4339200
$ws.Range("I58").Value = 3547500
$ws.Range("J58").Value = 2487300

# Row 59
$ws.Range("D59").Value = 839700
$ws.Range("E59").Value = 696200
$ws.Range("F59").Value = 466800
$ws.Range("G59").Value = 580100
$ws.Range("H59").Value = 504800
$ws.Range("I59").Value = 704400
$ws.Range("J59").Value = 491800

# Row 60
$ws.Range("D60").Value = 9350400
$ws.Range("E60").Value = 9457500
$ws.Range("F60").Value = 8057100
$ws.Range("G60").Value = 7625900
$ws.Range("H60").Value = 6744800
$ws.Range("I60").Value = 6241600
$ws.Range("J60").Value = 4642700

# Row 61
$ws.Range("D61").Value = 6103800
$ws.Range("E61").Value = 6355000
$ws.Range("F61").Value = 5033700
$ws.Range("G61").Value = 5170000
$ws.Range("H61").Value = 5988100
$ws.Range("I61").Value = 6857700
$ws.Range("J61").Value = 7090600

# Row 62
$ws.Range("D62").Value = 1522400
$ws.Range("E62").Value = 1643100
$ws.Range("F62").Value = 1339300
$ws.Range("G62").Value = 1162400
$ws.Range("H62").Value = 1166500
$ws.Range("I62").Value = 1247600
$ws.Range("J62").Value = 1144200

# Row 66
$ws.Range("D66").Value = 19873000
$ws.Range("E66").Value = 20366600
$ws.Range("F66").Value = 15966600
$ws.Range("G66").Value = 15382900
$ws.Range("H66").Value = 14763200
$ws.Range("I66").Value = 14630600
$ws.Range("J66").Value = 13221200

# Row 72
$ws.Range("D72").Value = 7195000
$ws.Range("E72").Value = 6827400
$ws.Range("F72").Value = 6171900
$ws.Range("G72").Value = 5747100
$ws.Range("H72").Value = 5555400
$ws.Range("I72").Value = 5506800
$ws.Range("J72").Value = 5477800

# Row 76
$ws.Range("D76").Value = 9257000
$ws.Range("E76").Value = 8916600
$ws.Range("F76").Value = 7750000
$ws.Range("G76").Value = 6500900
$ws.Range("H76").Value = 6369900
$ws.Range("I76").Value = 6513600
$ws.Range("J76").Value = 6724200

# Row 81
$ws.Range("D81").Value = 686100
$ws.Range("E81").Value = 580500
$ws.Range("F81").Value = 512200
$ws.Range("G81").Value = 189100
$ws.Range("H81").Value = 55700
$ws.Range("I81").Value = 31900
$ws.Range("J81").Value = -7500

# Row 89
$ws.Range("D89").Value = 1178600
$ws.Range("E89").Value = 1411400
$ws.Range("F89").Value = 393300
$ws.Range("G89").Value = 1676700
$ws.Range("H89").Value = -11500
$ws.Range("I89").Value = 604300
$ws.Range("J89").Value = -303100

# Row 91
$ws.Range("D91").Value = -70400
$ws.Range("E91").Value = 141800
$ws.Range("F91").Value = -276300
$ws.Range("G91").Value = 115500
$ws.Range("H91").Value = -365300
$ws.Range("I91").Value = -1506700
$ws.Range("J91").Value = -320200

# Row 94
$ws.Range("D94").Value = -175900
$ws.Range("E94").Value = -3544200
$ws.Range("F94").Value = -1032900
$ws.Range("G94").Value = -1469100
$ws.Range("H94").Value = -750400
$ws.Range("I94").Value = -658800
$ws.Range("J94").Value = -558700

# Row 96
$ws.Range("D96").Value = -356100
$ws.Range("E96").Value = -44500
$ws.Range("F96").Value = -37200
$ws.Range("H96").Value = -7900
$ws.Range("J96").Value = -14600

# Row 100
$ws.Range("D100").Value = -623400
$ws.Range("E100").Value = 2852700
$ws.Range("F100").Value = 646800
$ws.Range("G100").Value = 27500
$ws.Range("H100").Value = -21700
$ws.Range("I100").Value = 771000
$ws.Range("J100").Value = 872500

# Row 101
$ws.Range("D101").Value = 17400
$ws.Range("E101").Value = -79700
$ws.Range("F101").Value = 42800
$ws.Range("G101").Value = -4100
$ws.Range("I101").Value = 37200

# Row 102
$ws.Range("D102").Value = 396800
$ws.Range("E102").Value = 640200
$ws.Range("F102").Value = 50000
$ws.Range("G102").Value = 231000
$ws.Range("H102").Value = -787900
$ws.Range("I102").Value = 753700
$ws.Range("J102").Value = 8100
